$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new IPO record ("스튜디오삼익") was spliced into the demand-forecast
# table. Column A (stock names) keeps its row position; the date / price /
# amount / underwriter columns (B:F) for rows 4-13 roll down by one slot so
# the new record's own B:F data lands on row 4 (next to "IBKS스팩24호" in
# column A), and the previous row 13 data ("스튜디오삼익" 's own date bucket,
# 2024.01.05~01.11) is discarded off the bottom of that block.

$ws.Range("B4").Value = "2024.01.17~01.23"
$ws.Range("C4").Value = "14,500~16,500"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 12325
$ws.Range("F4").Value = "DB금융투자"

$ws.Range("B5").Value = "2024.01.17~01.18"
$ws.Range("C5").Value = "2,000~2,000"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 2147483647
$ws.Range("F5").Value = "아이비케이투자증권"

$ws.Range("B6").Value = "2024.01.12~01.18"
$ws.Range("C6").Value = "4,800~5,800"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 4800
$ws.Range("F6").Value = "NH투자증권"

$ws.Range("B7").Value = "2024.01.12~01.18"
$ws.Range("C7").Value = "12,000~14,000"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 7200
$ws.Range("F7").Value = "키움증권"

$ws.Range("B8").Value = "2024.01.11~01.17"
$ws.Range("C8").Value = "9,200~11,000"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = 27600
$ws.Range("F8").Value = "삼성증권"

$ws.Range("B9").Value = "2024.01.09~01.10"
$ws.Range("C9").Value = "2,000~2,000"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = 11000
$ws.Range("F9").Value = "대신증권"

$ws.Range("B10").Value = "2024.01.08~01.12"
$ws.Range("C10").Value = "2,400~2,800"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = 16000
$ws.Range("F10").Value = "NH투자증권"

$ws.Range("B11").Value = "2024.01.08~01.12"
$ws.Range("C11").Value = "5,000~6,300"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = 43535
$ws.Range("F11").Value = "미래에셋증권"

$ws.Range("B12").Value = "2024.01.08~01.12"
$ws.Range("C12").Value = "4,300~4,900"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = 8858
$ws.Range("F12").Value = "케이비증권"

$ws.Range("B13").Value = "2024.01.05~01.11"
$ws.Range("C13").Value = "13,000~15,000"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = 19500
$ws.Range("F13").Value = "하나증권"
